$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp column (O) for every data row (2 through 395)
# from "2023-01-25 06:49:22" to "2023-01-25 12:56:30"
$newTimestamp = "2023-01-25 12:56:30"
for ($row = 2; $row -le 395; $row++) {
    $ws.Range("O" + $row).Value = $newTimestamp
}

# Row 27: ratingAmount (D27) increased from 68 to 69
$ws.Range("D27").Value = 69

# Row 31: ratingAmount (D31) increased from 12 to 13
$ws.Range("D31").Value = 13

# Row 217: ratingAmount (D217) and ratingValue (E217) swapped (4/5 -> 5/4)
$ws.Range("D217").Value = 5
$ws.Range("E217").Value = 4

# Row 255: productAriaLabel (M255) updated with "Online kein Bestand" note
$ws.Range("M255").Value = "Pasquier Schokobrötchen 16 Stück - Online kein Bestand 8.50 Schweizer Franken"

# Row 281: ratingAmount (D281) increased from 1 to 2
$ws.Range("D281").Value = 2
